# Daily attendance processing - 2025-12-06 19:45:37
# Normalize the "Recorded By" (column G) entries so that "System" sorts
# ahead of named/email recorders, matching the updated attendance export
# ordering convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value()

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "backup@backdoor.com, system, System") {
        $cell.Value = "backup@backdoor.com, System, system"
    }
    elseif ($val -eq "dnasr281@gmail.com, admin@admin.com") {
        $cell.Value = "admin@admin.com, dnasr281@gmail.com"
    }
}
